$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.3604697925125322
$ws.Range("C2").Value = 0.05263273439551597
$ws.Range("E2").Value = 0.6815605889087237
$ws.Range("F2").Value = 2.278349838371255
$ws.Range("G2").Value = 0.002436492128489158
$ws.Range("I2").Value = 0.5262371975377107
$ws.Range("J2").Value = 0.04682457739595414
$ws.Range("K2").Value = 0.3771883825602629
$ws.Range("O2").Value = 2.400377726386964
# Row 3
$ws.Range("B3").Value = 0.3184191154466589
$ws.Range("C3").Value = 0.04602754622919747
$ws.Range("E3").Value = 0.6504940386193994
$ws.Range("F3").Value = 2.253495414606036
$ws.Range("G3").Value = 0.002438722578336879
$ws.Range("I3").Value = 0.534929670354833
$ws.Range("J3").Value = 0.04748271991215347
$ws.Range("K3").Value = 0.3312207902634725
$ws.Range("O3").Value = 2.433317797950323
# Row 4
$ws.Range("B4").Value = 0.2925570874589312
$ws.Range("C4").Value = 0.04195400596935883
$ws.Range("E4").Value = 0.6317205993841952
$ws.Range("F4").Value = 2.239648688146147
$ws.Range("G4").Value = 0.002440163241269938
$ws.Range("I4").Value = 0.5406116965563328
$ws.Range("J4").Value = 0.04791906164393822
$ws.Range("K4").Value = 0.3029298332953374
$ws.Range("O4").Value = 2.45518959098473
# Row 5
$ws.Range("B5").Value = 0.2820080124471929
$ws.Range("C5").Value = 0.04028957426199042
$ws.Range("E5").Value = 0.6241460730570907
$ws.Range("F5").Value = 2.234361074789476
$ws.Range("G5").Value = 0.00244076826975357
$ws.Range("I5").Value = 0.543013774307525
$ws.Range("J5").Value = 0.04810496889260207
$ws.Range("K5").Value = 0.2913849524832699
$ws.Range("O5").Value = 2.464515987183873
# Row 6
$ws.Range("B6").Value = 0.2802557567340784
$ws.Range("C6").Value = 0.04001293131065609
$ws.Range("E6").Value = 0.6228929064251645
$ws.Range("F6").Value = 2.233504499113266
$ws.Range("G6").Value = 0.002440869819663127
$ws.Range("I6").Value = 0.5434178653449315
$ws.Range("J6").Value = 0.04813632717105065
$ws.Range("K6").Value = 0.28946698231303
$ws.Range("O6").Value = 2.466089584156435
# Row 7
$ws.Range("B7").Value = 0.2924148589156346
$ws.Range("C7").Value = 0.0419315766987296
$ws.Range("E7").Value = 0.6316181397814802
$ws.Range("F7").Value = 2.23957594069185
$ws.Range("G7").Value = 0.00244017132808531
$ws.Range("I7").Value = 0.5406437412863276
$ws.Range("J7").Value = 0.04792153609687588
$ws.Range("K7").Value = 0.3027741992252686
$ws.Range("O7").Value = 2.455313696667346
# Row 8
$ws.Range("B8").Value = 0.3459799717795704
$ws.Range("C8").Value = 0.0503590306428805
$ws.Range("E8").Value = 0.6707862297667475
$ws.Range("F8").Value = 2.269486231788449
$ws.Range("G8").Value = 0.002437246451052265
$ws.Range("I8").Value = 0.5291627370127081
$ws.Range("J8").Value = 0.04704480615966311
$ws.Range("K8").Value = 0.3613529858972129
$ws.Range("O8").Value = 2.411393411253187
# Row 9
$ws.Range("B9").Value = 0.450658060906477
$ws.Range("C9").Value = 0.06674039478714633
$ws.Range("E9").Value = 0.7499932344588132
$ws.Range("F9").Value = 2.339389429036629
$ws.Range("G9").Value = 0.00243207293101233
$ws.Range("I9").Value = 0.5093886683923348
$ws.Range("J9").Value = 0.04558189953649894
$ws.Range("K9").Value = 0.475672619920033
$ws.Range("O9").Value = 2.338354189223651
# Row 10
$ws.Range("B10").Value = 0.5273194355882822
$ws.Range("C10").Value = 0.0786851612534889
$ws.Range("E10").Value = 0.809663779002733
$ws.Range("F10").Value = 2.397655674448714
$ws.Range("G10").Value = 0.002428611202677897
$ws.Range("I10").Value = 0.4965365244806375
$ws.Range("J10").Value = 0.04466413419630832
$ws.Range("K10").Value = 0.5593017620992669
$ws.Range("O10").Value = 2.292704443571367
# Row 11
$ws.Range("B11").Value = 0.5621366971575981
$ws.Range("C11").Value = 0.08409905619558344
$ws.Range("E11").Value = 0.8371338855577193
$ws.Range("F11").Value = 2.425674110578001
$ws.Range("G11").Value = 0.002427109304294464
$ws.Range("I11").Value = 0.4910548003638624
$ws.Range("J11").Value = 0.04428087916850032
$ws.Range("K11").Value = 0.5972637468115067
$ws.Range("O11").Value = 2.273684117483484
# Row 12
$ws.Range("B12").Value = 0.5753124256746105
$ws.Range("C12").Value = 0.08614624445311847
$ws.Range("E12").Value = 0.8475830860623006
$ws.Range("F12").Value = 2.436502261312256
$ws.Range("G12").Value = 0.002426550996170618
$ws.Range("I12").Value = 0.4890315784846688
$ws.Range("J12").Value = 0.04414068838239338
$ws.Range("K12").Value = 0.6116267031724476
$ws.Range("O12").Value = 2.266733305508239
# Row 13
$ws.Range("B13").Value = 0.5724751999288173
$ws.Range("C13").Value = 0.08570547774550619
$ws.Range("E13").Value = 0.8453305787390235
$ws.Range("F13").Value = 2.434160511679806
$ws.Range("G13").Value = 0.002426670774730604
$ws.Range("I13").Value = 0.4894649745260846
$ws.Range("J13").Value = 0.04417066103582279
$ws.Range("K13").Value = 0.6085339476104537
$ws.Range("O13").Value = 2.26821907835577
# Row 14
$ws.Range("B14").Value = 0.5632208529178797
$ws.Range("C14").Value = 0.08426753892752004
$ws.Range("E14").Value = 0.8379926080669406
$ws.Range("F14").Value = 2.426560572814509
$ws.Range("G14").Value = 0.002427063163147812
$ws.Range("I14").Value = 0.4908872939758275
$ws.Range("J14").Value = 0.04426924650092445
$ws.Range("K14").Value = 0.5984456490122056
$ws.Range("O14").Value = 2.273107220064119
# Row 15
$ws.Range("B15").Value = 0.5575511296173659
$ws.Range("C15").Value = 0.08338637503226209
$ws.Range("E15").Value = 0.8335039937647508
$ws.Range("F15").Value = 2.421933822805215
$ws.Range("G15").Value = 0.00242730486904743
$ws.Range("I15").Value = 0.4917653581803521
$ws.Range("J15").Value = 0.04433027669200662
$ws.Range("K15").Value = 0.5922646354088101
$ws.Range("O15").Value = 2.276134157147311
# Row 16
$ws.Range("B16").Value = 0.5250428507364404
$ws.Range("C16").Value = 0.07833094412508501
$ws.Range("E16").Value = 0.807875099341004
$ws.Range("F16").Value = 2.39585509983965
$ws.Range("G16").Value = 0.002428710816982548
$ws.Range("I16").Value = 0.4969021197313008
$ws.Range("J16").Value = 0.04468987117681777
$ws.Range("K16").Value = 0.55681915818829
$ws.Range("O16").Value = 2.29398266020732
# Row 17
$ws.Range("B17").Value = 0.5050851307044582
$ws.Range("C17").Value = 0.07522446014199602
$ws.Range("E17").Value = 0.7922360405728455
$ws.Range("F17").Value = 2.380244569510552
$ws.Range("G17").Value = 0.002429591946835595
$ws.Range("I17").Value = 0.5001468802336895
$ws.Range("J17").Value = 0.04491925188890988
$ws.Range("K17").Value = 0.535053155634813
$ws.Range("O17").Value = 2.305379824856928
# Row 18
$ws.Range("B18").Value = 0.4936007129160487
$ws.Range("C18").Value = 0.07343583029779666
$ws.Range("E18").Value = 0.7832715195239643
$ws.Range("F18").Value = 2.371408161751958
$ws.Range("G18").Value = 0.002430105610203414
$ws.Range("I18").Value = 0.5020475075570801
$ws.Range("J18").Value = 0.04505440759087342
$ws.Range("K18").Value = 0.5225263145787835
$ws.Range("O18").Value = 2.312099469537714
# Row 19
$ws.Range("B19").Value = 0.4897114036525352
$ws.Range("C19").Value = 0.07282991364552061
$ws.Range("E19").Value = 0.7802415501940487
$ws.Range("F19").Value = 2.368440736527262
$ws.Range("G19").Value = 0.002430280707629962
$ws.Range("I19").Value = 0.502696918941508
$ws.Range("J19").Value = 0.04510072188782388
$ws.Range("K19").Value = 0.5182836563122635
$ws.Range("O19").Value = 2.314402822505514
# Row 20
$ws.Range("B20").Value = 0.5072102141095343
$ws.Range("C20").Value = 0.07555534402470698
$ws.Range("E20").Value = 0.793897672472383
$ws.Range("F20").Value = 2.381891599174736
$ws.Range("G20").Value = 0.002429497439387039
$ws.Range("I20").Value = 0.4997979162688004
$ws.Range("J20").Value = 0.04489450034356324
$ws.Range("K20").Value = 0.5373709768733193
$ws.Range("O20").Value = 2.304149567993136
# Row 21
$ws.Range("B21").Value = 0.5659393237394852
$ws.Range("C21").Value = 0.08468997676521894
$ws.Range("E21").Value = 0.8401466756648972
$ws.Range("F21").Value = 2.428786932745254
$ws.Range("G21").Value = 0.002426947626660077
$ws.Range("I21").Value = 0.4904680959889447
$ws.Range("J21").Value = 0.04424015537593462
$ws.Range("K21").Value = 0.6014091698006894
$ws.Range("O21").Value = 2.271664615700942
# Row 22
$ws.Range("B22").Value = 0.6042704998586146
$ws.Range("C22").Value = 0.09064283432806519
$ws.Range("E22").Value = 0.8706462701416058
$ws.Range("F22").Value = 2.460707744587808
$ws.Range("G22").Value = 0.002425341938225488
$ws.Range("I22").Value = 0.4846771045691938
$ws.Range("J22").Value = 0.04384130209146697
$ws.Range("K22").Value = 0.6431890832253657
$ws.Range("O22").Value = 2.251901672723463
# Row 23
$ws.Range("B23").Value = 0.5838173914882248
$ws.Range("C23").Value = 0.08746727909073115
$ws.Range("E23").Value = 0.8543430524808713
$ws.Range("F23").Value = 2.443554407720143
$ws.Range("G23").Value = 0.002426193380604887
$ws.Range("I23").Value = 0.4877397690570593
$ws.Range("J23").Value = 0.04405153726657396
$ws.Range("K23").Value = 0.6208972529232142
$ws.Range("O23").Value = 2.26231498607649
# Row 24
$ws.Range("B24").Value = 0.5062494964266477
$ws.Range("C24").Value = 0.07540575973352759
$ws.Range("E24").Value = 0.7931463658947706
$ws.Range("F24").Value = 2.381146546326079
$ws.Range("G24").Value = 0.00242954014412735
$ws.Range("I24").Value = 0.4999555733292578
$ws.Range("J24").Value = 0.04490568029640229
$ws.Range("K24").Value = 0.5363231311914376
$ws.Range("O24").Value = 2.304705246151201
# Row 25
$ws.Range("B25").Value = 0.4223812540293466
$ws.Range("C25").Value = 0.06232455105326551
$ws.Range("E25").Value = 0.7283073284946084
$ws.Range("F25").Value = 2.319269248642087
$ws.Range("G25").Value = 0.002433412680350399
$ws.Range("I25").Value = 0.5144441696974056
$ws.Range("J25").Value = 0.04595013239689116
$ws.Range("K25").Value = 0.4448078775396027
$ws.Range("O25").Value = 2.356708613029625
